# #5: cash & deposit done
# Rework the "存款" (deposits) sheet: fix the header row (which had
# previously been a stray duplicate of the data row) to hold real column
# headers, and append the standard metadata columns (property_category,
# category, date, legislator_name, legislator_id, source_file, index)
# that the other sheets already carry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# ---- header row (row 1) ----
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# new header cells (G1:M1) need to look like the rest of row 1 (bold,
# thin box border, centered / top aligned) to match the B1:F1 style
$headerNew = $ws.Range("G1:M1")
$headerNew.Font.Bold = $true
$headerNew.Borders.LineStyle = 1
$headerNew.HorizontalAlignment = -4108
$headerNew.VerticalAlignment = -4160

# ---- row 2 (record 81) ----
$ws.Range("G2").Value = "deposit"
$ws.Range("H2").Value = "normal"
$ws.Range("I2").Value = "2011-11-23"
$ws.Range("J2").Value = "翁重鈞"
$ws.Range("K2").Value = 551
$ws.Range("L2").Value = "tmp22571"
$ws.Range("M2").Value = 81

# ---- row 3 (record 82) ----
$ws.Range("G3").Value = "deposit"
$ws.Range("H3").Value = "normal"
$ws.Range("I3").Value = "2011-11-23"
$ws.Range("J3").Value = "翁重鈞"
$ws.Range("K3").Value = 551
$ws.Range("L3").Value = "tmp22571"
$ws.Range("M3").Value = 82
